$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.144.72'
$ws.Range('E2').Value = '  +1.71%  '
$ws.Range('D3').Value = '3.215.58'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.04'
$ws.Range('E5').Value = '  +4.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.15'
$ws.Range('E6').Value = '  +1.89%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '3.215.48'
$ws.Range('E8').Value = '  +1.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.534'
$ws.Range('E9').Value = '  -0.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.18'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.509'
$ws.Range('E12').Value = '  +1.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000274'
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.63'
$ws.Range('E14').Value = '  +1.92%  '
$ws.Range('D15').Value = '3.737.05'
$ws.Range('E15').Value = '  +1.08%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.48'
$ws.Range('E16').Value = '  +4.15%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '66.214.94'
$ws.Range('E17').Value = '  +1.69%  '
$ws.Range('D18').Value = '3.211.82'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '512.77'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.59'
$ws.Range('E21').Value = '  +4.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.734'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.06'
$ws.Range('E23').Value = '  +3.22%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.28'
$ws.Range('E24').Value = '  +0.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.20'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('E27').Value = '  +2.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.22'
$ws.Range('E28').Value = '  +2.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.25'
$ws.Range('E29').Value = '  +2.83%  '
$ws.Range('B30').Value = 'Stacks'
$ws.Range('C30').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.88'
$ws.Range('E30').Value = '  +3.45%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.88'
$ws.Range('E31').Value = '  +8.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.20'
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('E33').Value = '  +0.91%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.62'
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.36'
$ws.Range('E36').Value = '  -0.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0911'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '483.18'
$ws.Range('E38').Value = '  +1.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0422'
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.98'
$ws.Range('E40').Value = '  -5.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.86'
$ws.Range('E41').Value = '  +2.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.299'
$ws.Range('E42').Value = '  +4.10%  '
$ws.Range('E43').Value = '  -0.49%  '
$ws.Range('D44').Value = '2.946.37'
$ws.Range('E44').Value = '  -3.96%  '
$ws.Range('E45').Value = '  +1.72%  '
$ws.Range('D46').Value = '0.0₃0643'
$ws.Range('E46').Value = '  +4.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.85'
$ws.Range('E47').Value = '  -1.06%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.116'
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('E50').Value = '  +2.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.56'
$ws.Range('E51').Value = '  -1.03%  '
